$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "09.25.2019"
$ws.Range("C1").Value = "10.07.2019"
$ws.Range("D1").Value = "delta(12 Days)"
$ws.Range("E1").Value = "2019.10.11"
$ws.Range("F1").Value = "delat (16 Days)"

$ws.Range("A2").Value = "歪歪歪"
$ws.Range("B2").Value = 222546
$ws.Range("C2").Value = 227745
$ws.Range("D2").Formula = "=C2-B2"
$ws.Range("E2").Value = 229162
$ws.Range("F2").Formula = "=E2-B2"

$ws.Range("A3").Value = "661工作室的小刘"
$ws.Range("B3").Value = 191385
$ws.Range("C3").Value = 195531
$ws.Range("D3").Formula = "=C3-B3"
$ws.Range("E3").Value = 197259
$ws.Range("F3").Formula = "=E3-B3"

$ws.Range("A4").Value = "柴犬婧"
$ws.Range("B4").Value = 156090
$ws.Range("C4").Value = 162090
$ws.Range("D4").Formula = "=C4-B4"
$ws.Range("E4").Value = 164090
$ws.Range("F4").Formula = "=E4-B4"

$ws.Range("A5").Value = "feifei2"
$ws.Range("B5").Value = 159453
$ws.Range("C5").Value = 162343
$ws.Range("D5").Formula = "=C5-B5"
$ws.Range("E5").Value = 163794
$ws.Range("F5").Formula = "=E5-B5"

$ws.Range("A6").Value = "艳玲"
$ws.Range("B6").Value = 154882
$ws.Range("C6").Value = 159902
$ws.Range("D6").Formula = "=C6-B6"
$ws.Range("E6").Value = 161702
$ws.Range("F6").Formula = "=E6-B6"

$ws.Range("A7").Value = "李同学"
$ws.Range("B7").Value = 125566
$ws.Range("C7").Value = 129799
$ws.Range("D7").Formula = "=C7-B7"
$ws.Range("E7").Value = 130827
$ws.Range("F7").Formula = "=E7-B7"

$ws.Range("A8").Value = "万欢"
$ws.Range("B8").Value = 107212
$ws.Range("C8").Value = 110348
$ws.Range("D8").Formula = "=C8-B8"
$ws.Range("E8").Value = 111000
$ws.Range("F8").Formula = "=E8-B8"

$ws.Range("A9").Value = "艳玲"
$ws.Range("B9").Value = 97112
$ws.Range("C9").Value = 99946
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("E9").Value = 101195
$ws.Range("F9").Formula = "=E9-B9"

$ws.Range("A10").Value = "快乐的蓝丁丁"
$ws.Range("B10").Value = 93170
$ws.Range("C10").Value = 95670
$ws.Range("D10").Formula = "=C10-B10"
$ws.Range("E10").Value = 96270
$ws.Range("F10").Formula = "=E10-B10"

$ws.Range("A11").Value = "和上帝旅行"
$ws.Range("B11").Value = 85328
$ws.Range("C11").Value = 91328
$ws.Range("D11").Formula = "=C11-B11"
$ws.Range("E11").Value = 93328
$ws.Range("F11").Formula = "=E11-B11"

$ws.Range("A12").Value = "叶子"
$ws.Range("B12").Value = 85800
$ws.Range("C12").Value = 88300
$ws.Range("D12").Formula = "=C12-B12"
$ws.Range("E12").Value = 88800
$ws.Range("F12").Formula = "=E12-B12"

$ws.Range("A13").Value = "冷了蜘蛛"
$ws.Range("B13").Value = 85387
$ws.Range("C13").Value = 87860
$ws.Range("D13").Formula = "=C13-B13"
$ws.Range("E13").Value = 88521
$ws.Range("F13").Formula = "=E13-B13"

$ws.Range("A14").Value = "小天才"
$ws.Range("B14").Value = 83037
$ws.Range("C14").Value = 84650
$ws.Range("D14").Formula = "=C14-B14"
$ws.Range("E14").Value = 84906
$ws.Range("F14").Formula = "=E14-B14"

$ws.Range("A15").Value = "阿牙布"
$ws.Range("B15").Value = 76033
$ws.Range("C15").Value = 77233
$ws.Range("D15").Formula = "=C15-B15"
$ws.Range("E15").Value = 77333
$ws.Range("F15").Formula = "=E15-B15"

$ws.Range("A16").Value = "其乐"
$ws.Range("B16").Value = 70236
$ws.Range("C16").Value = 71342
$ws.Range("D16").Formula = "=C16-B16"
$ws.Range("E16").Value = 71687
$ws.Range("F16").Formula = "=E16-B16"

$ws.Range("A17").Value = "阿黄的马甲"
$ws.Range("B17").Value = 63345
$ws.Range("C17").Value = 66310
$ws.Range("D17").Formula = "=C17-B17"
$ws.Range("E17").Value = 66810
$ws.Range("F17").Formula = "=E17-B17"

$ws.Range("A18").Value = "芳芳"
$ws.Range("B18").Value = 62100
$ws.Range("C18").Value = 64500
$ws.Range("D18").Formula = "=C18-B18"
$ws.Range("E18").Value = 64900
$ws.Range("F18").Formula = "=E18-B18"

$ws.Range("A19").Value = "popovicki"
$ws.Range("B19").Value = 57102
$ws.Range("C19").Value = 59607
$ws.Range("D19").Formula = "=C19-B19"
$ws.Range("E19").Value = 60572
$ws.Range("F19").Formula = "=E19-B19"

$ws.Range("A20").Value = "小戴"
$ws.Range("B20").Value = 59518
$ws.Range("C20").Value = 60225
$ws.Range("D20").Formula = "=C20-B20"
$ws.Range("E20").Value = 60511
$ws.Range("F20").Formula = "=E20-B20"

$ws.Range("A21").Value = "除除草"
$ws.Range("B21").Value = 52663
$ws.Range("C21").Value = 55543
$ws.Range("D21").Formula = "=C21-B21"
$ws.Range("E21").Value = 56588
$ws.Range("F21").Formula = "=E21-B21"

$ws.Range("A22").Value = "角"
$ws.Range("B22").Value = 54233
$ws.Range("C22").Value = 56033
$ws.Range("D22").Formula = "=C22-B22"
$ws.Range("E22").Value = 56533
$ws.Range("F22").Formula = "=E22-B22"

$ws.Range("A23").Value = "弘"
$ws.Range("B23").Value = 51280
$ws.Range("C23").Value = 53080
$ws.Range("D23").Formula = "=C23-B23"
$ws.Range("E23").Value = 53530
$ws.Range("F23").Formula = "=E23-B23"

$ws.Range("A24").Value = "诺瑶"
$ws.Range("B24").Value = 49473
$ws.Range("C24").Value = 50973
$ws.Range("D24").Formula = "=C24-B24"
$ws.Range("E24").Value = 51973
$ws.Range("F24").Formula = "=E24-B24"

$ws.Range("A25").Value = "鱼姑娘"
$ws.Range("B25").Value = 48200
$ws.Range("C25").Value = 49700
$ws.Range("D25").Formula = "=C25-B25"
$ws.Range("E25").Value = 50200
$ws.Range("F25").Formula = "=E25-B25"

$ws.Range("A26").Value = "leeRoc"
$ws.Range("B26").Value = 42530
$ws.Range("C26").Value = 43730
$ws.Range("D26").Formula = "=C26-B26"
$ws.Range("E26").Value = 44130
$ws.Range("F26").Formula = "=E26-B26"

$ws.Range("A27").Value = "Fish"
$ws.Range("B27").Value = 40671
$ws.Range("C27").Value = 42271
$ws.Range("D27").Formula = "=C27-B27"
$ws.Range("E27").Value = 43171
$ws.Range("F27").Formula = "=E27-B27"

$ws.Range("A28").Value = "顺风逆风"
$ws.Range("B28").Value = 36340
$ws.Range("C28").Value = 38840
$ws.Range("D28").Formula = "=C28-B28"
$ws.Range("E28").Value = 39340
$ws.Range("F28").Formula = "=E28-B28"

$ws.Range("A29").Value = "飞扬"
$ws.Range("B29").Value = 34448
$ws.Range("C29").Value = 36848
$ws.Range("D29").Formula = "=C29-B29"
$ws.Range("E29").Value = 37648
$ws.Range("F29").Formula = "=E29-B29"

$ws.Range("A30").Value = "安静的沼泽"
$ws.Range("B30").Value = 32230
$ws.Range("C30").Value = 33880
$ws.Range("D30").Formula = "=C30-B30"
$ws.Range("E30").Value = 34330
$ws.Range("F30").Formula = "=E30-B30"

$ws.Range("A31").Value = "风风"
$ws.Range("B31").Value = 31008
$ws.Range("C31").Value = 32408
$ws.Range("D31").Formula = "=C31-B31"
$ws.Range("E31").Value = 33008
$ws.Range("F31").Formula = "=E31-B31"

$ws.Range("A32").Value = "精灵海"
$ws.Range("B32").Value = 27100
$ws.Range("C32").Value = 28600
$ws.Range("D32").Formula = "=C32-B32"
$ws.Range("E32").Value = 29100
$ws.Range("F32").Formula = "=E32-B32"

$ws.Range("A33").Value = "小兰"
$ws.Range("B33").Value = 26550
$ws.Range("C33").Value = 28150
$ws.Range("D33").Formula = "=C33-B33"
$ws.Range("E33").Value = 28650
$ws.Range("F33").Formula = "=E33-B33"

$ws.Range("A34").Value = "线"
$ws.Range("B34").Value = 15601
$ws.Range("C34").Value = 15898
$ws.Range("D34").Formula = "=C34-B34"
$ws.Range("E34").Value = 16000
$ws.Range("F34").Formula = "=E34-B34"

$ws.Range("A35").Value = "lyl"
$ws.Range("B35").Value = "N/A"
$ws.Range("C35").Value = 13100
$ws.Range("E35").Value = 13600

$ws.Range("A36").Value = "婷婷"
$ws.Range("B36").Value = 10960
$ws.Range("C36").Value = 12360
$ws.Range("D36").Formula = "=C36-B36"
$ws.Range("E36").Value = 12960
$ws.Range("F36").Formula = "=E36-B36"

$ws.Range("A37").Value = "萌*6哥"
$ws.Range("B37").Value = 510
$ws.Range("C37").Value = 510
$ws.Range("D37").Formula = "=C37-B37"
$ws.Range("E37").Value = 510
$ws.Range("F37").Formula = "=E37-B37"

